# Adds the "ODI Bowling Extra" worksheet (new sheetId=5) to the workbook,
# as the last/trailing sheet after "ODI Batting Extra", and populates it
# with MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new worksheet right after the current last sheet
#    ("ODI Batting Extra"), then rename it.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# ---------------------------------------------------------------------
# 2. Header row (bold, matching the style used by the other "*Extra"
#    sheets in this workbook).
# ---------------------------------------------------------------------
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.NumberFormat = "@"
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# ---------------------------------------------------------------------
# 3. Data rows 2-21 -- every value in this sheet is stored as text
#    (match codes, maiden-over counts, and percentages alike), so force
#    a text number-format before writing each value.
# ---------------------------------------------------------------------
$data = @(
    @("3785", "0", ""),
    @("3843", "0", ""),
    @("3846", "0", "20.00%"),
    @("3847", "", ""),
    @("3849", "", ""),
    @("3850", "0", ""),
    @("3864", "0", ""),
    @("3867", "0", ""),
    @("3869", "0", ""),
    @("3916", "", ""),
    @("3941", "", ""),
    @("4006", "0", ""),
    @("4007", "0", ""),
    @("4008", "", ""),
    @("4010", "0", ""),
    @("4040", "0", ""),
    @("4198", "0", ""),
    @("4202", "0", ""),
    @("4262", "", ""),
    @("4340", "0", "")
)

$row = 2
foreach ($rec in $data) {
    $rowRange = $ws.Range("A" + $row + ":C" + $row)
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $rec[0]
    if ($rec[1] -ne "") {
        $ws.Cells.Item($row, 2).Value = $rec[1]
    }
    if ($rec[2] -ne "") {
        $ws.Cells.Item($row, 3).Value = $rec[2]
    }

    $row = $row + 1
}

# ---------------------------------------------------------------------
# 4. Leave the workbook selection on A1 of the new sheet, same as the
#    other sheets in the file.
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
